$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.819045543670654
$ws.Range("B1").Value = 3.486467123031616
$ws.Range("C1").Value = 1.883810043334961
$ws.Range("D1").Value = 1.500928044319153
$ws.Range("E1").Value = 1.38044261932373
